$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 666. This shifts the existing rows
# 666..739 down to 667..740, growing the used range to A1:R740.
$ws.Rows(666).Insert()

# Populate the newly inserted row 666 with its data.
$ws.Cells.Item(666, 1).Value  = 5
$ws.Cells.Item(666, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(666, 3).Value  = "Maule"
$ws.Cells.Item(666, 4).Value  = 45194
$ws.Cells.Item(666, 5).Value  = 7
$ws.Cells.Item(666, 6).Value  = 100112043
$ws.Cells.Item(666, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(666, 8).Value  = "Sin especificar"
$ws.Cells.Item(666, 9).Value  = "Primera"
$ws.Cells.Item(666, 10).Value = 500
$ws.Cells.Item(666, 11).Value = 10000
$ws.Cells.Item(666, 12).Value = 10000
$ws.Cells.Item(666, 13).Value = 10000
$ws.Cells.Item(666, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(666, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(666, 16).Value = 167
$ws.Cells.Item(666, 17).Value = 60
$ws.Cells.Item(666, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used
# throughout column D.
$ws.Cells.Item(666, 4).NumberFormat = $ws.Cells.Item(667, 4).NumberFormat
